$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.932.75"
$ws.Range("E2").Value = "  -3.26%  "

# Row 3
$ws.Range("D3").Value = "2.919.25"
$ws.Range("E3").Value = "  -3.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'585.70"
$ws.Range("E5").Value = "  -1.44%  "

# Row 6
$ws.Range("D6").Value = "'145.43"
$ws.Range("E6").Value = "  -5.72%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -2.31%  "

# Row 9
$ws.Range("D9").Value = "2.917.74"
$ws.Range("E9").Value = "  -3.93%  "

# Row 10
$ws.Range("D10").Value = "'6.87"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("D11").Value = "'0.144"
$ws.Range("E11").Value = "  -4.76%  "

# Row 12
$ws.Range("D12").Value = "'0.447"
$ws.Range("E12").Value = "  -3.88%  "

# Row 13
$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = "  -3.93%  "

# Row 14
$ws.Range("D14").Value = "'33.58"
$ws.Range("E14").Value = "  -6.01%  "

# Row 15
$ws.Range("E15").Value = "  +0.17%  "

# Row 16
$ws.Range("D16").Value = "3.404.26"
$ws.Range("E16").Value = "  -3.83%  "

# Row 17
$ws.Range("D17").Value = "60.905.59"
$ws.Range("E17").Value = "  -3.23%  "

# Row 18
$ws.Range("D18").Value = "'6.76"
$ws.Range("E18").Value = "  -4.43%  "

# Row 19
$ws.Range("D19").Value = "2.921.94"
$ws.Range("E19").Value = "  -3.76%  "

# Row 20
$ws.Range("D20").Value = "'430.44"
$ws.Range("E20").Value = "  -5.32%  "

# Row 21
$ws.Range("D21").Value = "'13.60"
$ws.Range("E21").Value = "  -4.76%  "

# Row 22
$ws.Range("D22").Value = "'0.680"
$ws.Range("E22").Value = "  -2.64%  "

# Row 23
$ws.Range("D23").Value = "'7.13"
$ws.Range("E23").Value = "  -5.37%  "

# Row 24
$ws.Range("D24").Value = "'80.48"
$ws.Range("E24").Value = "  -3.24%  "

# Row 25
$ws.Range("D25").Value = "'10.81"
$ws.Range("E25").Value = "  -3.91%  "

# Row 26
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  -3.22%  "

# Row 27
$ws.Range("D27").Value = "'11.95"
$ws.Range("E27").Value = "  -3.36%  "

# Row 28
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("E29").Value = "  +0.13%  "

# Row 30
$ws.Range("D30").Value = "'7.18"
$ws.Range("E30").Value = "  -4.57%  "

# Row 31
$ws.Range("D31").Value = "'2.61"
$ws.Range("E31").Value = "  -3.38%  "

# Row 32
$ws.Range("E32").Value = "  -3.39%  "

# Row 33
$ws.Range("D33").Value = "'26.55"
$ws.Range("E33").Value = "  -3.90%  "

# Row 34
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -3.73%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0867"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("E36").Value = "  -2.99%  "

# Row 37
$ws.Range("D37").Value = "'5.65"
$ws.Range("E37").Value = "  -4.97%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.127"
$ws.Range("E38").Value = "  -3.32%  "

# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  -6.48%  "

# Row 40
$ws.Range("D40").Value = "'49.56"
$ws.Range("E40").Value = "  -1.72%  "

# Row 41
$ws.Range("E41").Value = "  -5.40%  "

# Row 42
$ws.Range("D42").Value = "'8.65"
$ws.Range("E42").Value = "  -5.08%  "

# Row 43
$ws.Range("D43").Value = "'0.295"
$ws.Range("E43").Value = "  -2.53%  "

# Row 44
$ws.Range("D44").Value = "'41.51"
$ws.Range("E44").Value = "  -2.86%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0350"
$ws.Range("E45").Value = "  -3.04%  "

# Row 46
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'377.38"
$ws.Range("E46").Value = "  -4.42%  "

# Row 47
$ws.Range("D47").Value = "2.700.44"
$ws.Range("E47").Value = "  -0.94%  "

# Row 48
$ws.Range("D48").Value = "'132.69"
$ws.Range("E48").Value = "  +0.20%  "

# Row 50
$ws.Range("D50").Value = "'24.86"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  -2.36%  "
